$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.461.24'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.579.17'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '207.53'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -1.49%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.21'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0589'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '1.803.74'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '1.601.87'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('D16').Value = '27.460.03'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.93'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '214.27'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').Value = '0.0₃0689'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -2.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.73'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.87'
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.89'
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.04'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').Value = '1.362.33'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.529'
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.972'
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '64.05'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.76'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.16'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('E46').Value = '  -2.21%  '
$ws.Range('D47').Value = '1.715.85'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.31'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').Value = '0.0₇0996'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('E51').Value = '  -0.61%  '
